$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.888.12"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "2.046.20"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("E6").Value = "  -1.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.67"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.33%  "

$ws.Range("E10").Value = "  -1.29%  "

$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "2.348.63"
$ws.Range("E12").Value = "  -0.66%  "

$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.781"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("D17").Value = "2.065.26"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "37.846.90"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("E20").Value = "  -4.24%  "

$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.23%  "

$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("E30").Value = "  -2.30%  "

$ws.Range("E31").Value = "  -0.60%  "

$ws.Range("E32").Value = "  +8.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0593"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("E36").Value = "  +5.05%  "

$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.32%  "

$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.92%  "

$ws.Range("D41").Value = "1.537.21"
$ws.Range("E41").Value = "  +0.69%  "

$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.56%  "

$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0913"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.71%  "

$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("D51").Value = "2.237.45"
$ws.Range("E51").Value = "  -0.65%  "
